# Update the multiplication problems throughout the worksheet to the
# regenerated set of operands (commit: "Update master to output
# generated at c986bee"). Each "old=" expression appears exactly once in
# the document, so a sequence of single-shot Find/Replace calls on the
# whole document Range is unambiguous and order-independent.

$d = $word.ActiveDocument

$d.Content.Find.Execute("921×4=", $false, $false, $false, $false, $false, $true, 1, $false, "115×9=", 2)
$d.Content.Find.Execute("289×3=", $false, $false, $false, $false, $false, $true, 1, $false, "966×7=", 2)
$d.Content.Find.Execute("148×9=", $false, $false, $false, $false, $false, $true, 1, $false, "136×5=", 2)
$d.Content.Find.Execute("152×2=", $false, $false, $false, $false, $false, $true, 1, $false, "232×7=", 2)
$d.Content.Find.Execute("956×2=", $false, $false, $false, $false, $false, $true, 1, $false, "326×3=", 2)
$d.Content.Find.Execute("691×7=", $false, $false, $false, $false, $false, $true, 1, $false, "856×5=", 2)
$d.Content.Find.Execute("601×4=", $false, $false, $false, $false, $false, $true, 1, $false, "471×3=", 2)
$d.Content.Find.Execute("267×4=", $false, $false, $false, $false, $false, $true, 1, $false, "398×2=", 2)
$d.Content.Find.Execute("665×3=", $false, $false, $false, $false, $false, $true, 1, $false, "662×8=", 2)
$d.Content.Find.Execute("677×9=", $false, $false, $false, $false, $false, $true, 1, $false, "103×6=", 2)
$d.Content.Find.Execute("647×7=", $false, $false, $false, $false, $false, $true, 1, $false, "169×6=", 2)
$d.Content.Find.Execute("990×4=", $false, $false, $false, $false, $false, $true, 1, $false, "398×8=", 2)
$d.Content.Find.Execute("221×9=", $false, $false, $false, $false, $false, $true, 1, $false, "302×5=", 2)
$d.Content.Find.Execute("539×2=", $false, $false, $false, $false, $false, $true, 1, $false, "800×4=", 2)
$d.Content.Find.Execute("709×8=", $false, $false, $false, $false, $false, $true, 1, $false, "911×3=", 2)
$d.Content.Find.Execute("383×9=", $false, $false, $false, $false, $false, $true, 1, $false, "453×7=", 2)
$d.Content.Find.Execute("638×3=", $false, $false, $false, $false, $false, $true, 1, $false, "471×8=", 2)
$d.Content.Find.Execute("754×6=", $false, $false, $false, $false, $false, $true, 1, $false, "665×7=", 2)
$d.Content.Find.Execute("132×8=", $false, $false, $false, $false, $false, $true, 1, $false, "887×2=", 2)
$d.Content.Find.Execute("922×5=", $false, $false, $false, $false, $false, $true, 1, $false, "811×5=", 2)
$d.Content.Find.Execute("221×2=", $false, $false, $false, $false, $false, $true, 1, $false, "382×2=", 2)
$d.Content.Find.Execute("741×5=", $false, $false, $false, $false, $false, $true, 1, $false, "238×3=", 2)
$d.Content.Find.Execute("137×3=", $false, $false, $false, $false, $false, $true, 1, $false, "499×2=", 2)
$d.Content.Find.Execute("355×6=", $false, $false, $false, $false, $false, $true, 1, $false, "958×4=", 2)
$d.Content.Find.Execute("811×6=", $false, $false, $false, $false, $false, $true, 1, $false, "722×5=", 2)
